# Auto-generated edit script: refreshes market-price-derived profit columns (H:N)
# across the Pandaemonium_Profits leve-crafting sheets, per scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 4047.516
$ws.Range("I80").Value = 316.46155
$ws.Range("J80").Value = 6742.1665
$ws.Range("K80").Value = 949.38465
$ws.Range("L80").Value = 20226.4995
$ws.Range("M80").Value = 48.61535000000003
$ws.Range("N80").Value = -22222.4995
# Row 83
$ws.Range("H83").Value = 4047.516
$ws.Range("I83").Value = 316.46155
$ws.Range("J83").Value = 6742.1665
$ws.Range("K83").Value = 2848.15395
$ws.Range("L83").Value = 60679.4985
$ws.Range("M83").Value = 2143.84605
$ws.Range("N83").Value = -70663.4985
# Row 105
$ws.Range("H105").Value = 38935.5
$ws.Range("J105").Value = 38935.5
$ws.Range("L105").Value = 38935.5
$ws.Range("N105").Value = -45923.5
# Row 129
$ws.Range("H129").Value = 1099.6389
$ws.Range("J129").Value = 1132.9565
$ws.Range("L129").Value = 3398.8695
$ws.Range("N129").Value = -13398.8695
# Row 133
$ws.Range("H133").Value = 59999.375
$ws.Range("J133").Value = 59999.375
$ws.Range("L133").Value = 59999.375
$ws.Range("N133").Value = -70119.375
# Row 137
$ws.Range("H137").Value = 1701.54
$ws.Range("I137").Value = 1386.4615
$ws.Range("K137").Value = 4159.3845
$ws.Range("M137").Value = -1609.3845

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1679.5
$ws.Range("I2").Value = 1894.8889
$ws.Range("J2").Value = 1033.3334
$ws.Range("K2").Value = 1894.8889
$ws.Range("L2").Value = 1033.3334
$ws.Range("M2").Value = -1781.8889
$ws.Range("N2").Value = -1259.3334
# Row 32
$ws.Range("H32").Value = 20574.139
$ws.Range("I32").Value = 22711.576
$ws.Range("K32").Value = 22711.576
$ws.Range("M32").Value = -22424.576
# Row 74
$ws.Range("H74").Value = 4745.1177
$ws.Range("I74").Value = 1837.4231
$ws.Range("J74").Value = 14195.125
$ws.Range("K74").Value = 1837.4231
$ws.Range("L74").Value = 14195.125
$ws.Range("M74").Value = -963.4231
$ws.Range("N74").Value = -15943.125
# Row 77
$ws.Range("H77").Value = 4745.1177
$ws.Range("I77").Value = 1837.4231
$ws.Range("J77").Value = 14195.125
$ws.Range("K77").Value = 9187.1155
$ws.Range("L77").Value = 70975.625
$ws.Range("M77").Value = -4819.1155
$ws.Range("N77").Value = -79711.625
# Row 116
$ws.Range("H116").Value = 1679.5
$ws.Range("I116").Value = 1894.8889
$ws.Range("J116").Value = 1033.3334
$ws.Range("K116").Value = 1894.8889
$ws.Range("L116").Value = 1033.3334
$ws.Range("M116").Value = 399.1111000000001
$ws.Range("N116").Value = -5621.3334
# Row 128
$ws.Range("H128").Value = 75000
$ws.Range("J128").Value = 75000
$ws.Range("L128").Value = 75000
$ws.Range("N128").Value = -84960

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1679.5
$ws.Range("I3").Value = 1894.8889
$ws.Range("J3").Value = 1033.3334
$ws.Range("K3").Value = 1894.8889
$ws.Range("L3").Value = 1033.3334
$ws.Range("M3").Value = -1780.8889
$ws.Range("N3").Value = -1261.3334
# Row 75
$ws.Range("H75").Value = 19405.834
$ws.Range("I75").Value = 2199
$ws.Range("J75").Value = 22847.2
$ws.Range("K75").Value = 2199
$ws.Range("L75").Value = 22847.2
$ws.Range("M75").Value = -1263
$ws.Range("N75").Value = -24719.2
# Row 78
$ws.Range("H78").Value = 19405.834
$ws.Range("I78").Value = 2199
$ws.Range("J78").Value = 22847.2
$ws.Range("K78").Value = 6597
$ws.Range("L78").Value = 68541.60000000001
$ws.Range("M78").Value = -1917
$ws.Range("N78").Value = -77901.60000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2147.3455
$ws.Range("I31").Value = 1587.4166
$ws.Range("K31").Value = 1587.4166
$ws.Range("M31").Value = -1292.4166
# Row 34
$ws.Range("H34").Value = 2147.3455
$ws.Range("I34").Value = 1587.4166
$ws.Range("K34").Value = 1587.4166
$ws.Range("M34").Value = -1385.4166
# Row 69
$ws.Range("H69").Value = 18091
$ws.Range("I69").Value = 6182
$ws.Range("J69").Value = 30000
$ws.Range("K69").Value = 6182
$ws.Range("L69").Value = 30000
$ws.Range("M69").Value = -5433
$ws.Range("N69").Value = -31498
# Row 72
$ws.Range("H72").Value = 18091
$ws.Range("I72").Value = 6182
$ws.Range("J72").Value = 30000
$ws.Range("K72").Value = 18546
$ws.Range("L72").Value = 90000
$ws.Range("M72").Value = -14802
$ws.Range("N72").Value = -97488
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# Row 132
$ws.Range("H132").Value = 3800.0908
$ws.Range("I132").Value = 4162.237
$ws.Range("J132").Value = 2990.5881
$ws.Range("K132").Value = 12486.711
$ws.Range("L132").Value = 8971.764299999999
$ws.Range("M132").Value = -9956.710999999999
$ws.Range("N132").Value = -14031.7643
# Row 134
$ws.Range("H134").Value = 2332.2622
$ws.Range("I134").Value = 1452.1111
$ws.Range("J134").Value = 3599.68
$ws.Range("K134").Value = 4356.3333
$ws.Range("L134").Value = 10799.04
$ws.Range("M134").Value = -1821.3333
$ws.Range("N134").Value = -15869.04

$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 1800
$ws.Range("J25").Value = 5000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15338
# Row 30
$ws.Range("H30").Value = 1800
$ws.Range("J30").Value = 5000
$ws.Range("L30").Value = 15000
$ws.Range("N30").Value = -15204
# Row 40
$ws.Range("H40").Value = 64.5
# Row 46
$ws.Range("H46").Value = 2994.5454
$ws.Range("J46").Value = 2994.5454
$ws.Range("L46").Value = 8983.636200000001
$ws.Range("N46").Value = -9165.636200000001
# Row 58
$ws.Range("H58").Value = 3107.4075
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3107.4075
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 9322.2225
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -9578.2225
# Row 117
$ws.Range("H117").Value = 1482
$ws.Range("J117").Value = 1966.3334
$ws.Range("L117").Value = 5899.0002
$ws.Range("N117").Value = -12783.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 14
$ws.Range("H14").Value = 178336670
$ws.Range("I14").Value = 178336670
$ws.Range("K14").Value = 178336670
$ws.Range("M14").Value = -178336502
# Row 123
$ws.Range("H123").Value = 10326
$ws.Range("J123").Value = 10326
$ws.Range("L123").Value = 10326
$ws.Range("N123").Value = -15226

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 9336.333000000001
$ws.Range("J3").Value = 13002.5
$ws.Range("L3").Value = 13002.5
$ws.Range("N3").Value = -13226.5
# Row 15
$ws.Range("H15").Value = 9336.333000000001
$ws.Range("J15").Value = 13002.5
$ws.Range("L15").Value = 13002.5
$ws.Range("N15").Value = -13342.5
# Row 122
$ws.Range("H122").Value = 6080.817
$ws.Range("I122").Value = 5730.479
$ws.Range("J122").Value = 6811.9565
$ws.Range("K122").Value = 17191.437
$ws.Range("L122").Value = 20435.8695
$ws.Range("M122").Value = -14741.437
$ws.Range("N122").Value = -25335.8695
# Row 132
$ws.Range("H132").Value = 5362.811
$ws.Range("I132").Value = 6046.2856
$ws.Range("J132").Value = 4465.75
$ws.Range("K132").Value = 18138.8568
$ws.Range("L132").Value = 13397.25
$ws.Range("M132").Value = -15608.8568
$ws.Range("N132").Value = -18457.25
# Row 136
$ws.Range("H136").Value = 4236.9434
$ws.Range("I136").Value = 2517.8215
$ws.Range("J136").Value = 6162.36
$ws.Range("K136").Value = 7553.4645
$ws.Range("L136").Value = 18487.08
$ws.Range("M136").Value = -5003.4645
$ws.Range("N136").Value = -23587.08

$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 4840.5
$ws.Range("I136").Value = 3930.6924
$ws.Range("J136").Value = 5985.0967
$ws.Range("K136").Value = 11792.0772
$ws.Range("L136").Value = 17955.2901
$ws.Range("M136").Value = -9242.0772
$ws.Range("N136").Value = -23055.2901
